# Updates cryptos list price/volume data (GitHub Actions scrape refresh).
# Matches the per-cell text updates from the commit diff; B/C swaps on some
# rows reorder a few coins that changed rank between scrapes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.861.53'
$ws.Range("E2").Value = '  -3.84%  '
$ws.Range("D3").Value = '1.953.83'
$ws.Range("E3").Value = '  -3.90%  '
$ws.Range("D4").Value = '''1.01'
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").Value = '''240.88'
$ws.Range("E5").Value = '  -4.53%  '
$ws.Range("D6").Value = '''0.616'
$ws.Range("E6").Value = '  -4.45%  '
$ws.Range("D7").Value = '''59.98'
$ws.Range("E7").Value = '  -5.49%  '
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").Value = '''0.366'
$ws.Range("E9").Value = '  -1.84%  '
$ws.Range("D10").Value = '''56.10'
$ws.Range("E10").Value = '  -4.61%  '
$ws.Range("D11").Value = '''0.0780'
$ws.Range("E11").Value = '  +3.90%  '
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").Value = '''0.849'
$ws.Range("E13").Value = '  -6.05%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''13.71'
$ws.Range("E14").Value = '  -8.08%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.257.92'
$ws.Range("E15").Value = '  -3.30%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '''21.37'
$ws.Range("E16").Value = '  +4.73%  '
$ws.Range("D17").Value = '''5.33'
$ws.Range("E17").Value = '  -3.83%  '
$ws.Range("D18").Value = '1.979.94'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = '35.834.45'
$ws.Range("E19").Value = '  -3.82%  '
$ws.Range("D20").Value = '''70.22'
$ws.Range("E20").Value = '  -3.97%  '
$ws.Range("D21").Value = '0.0₃0838'
$ws.Range("E21").Value = '  -3.65%  '
$ws.Range("D22").Value = '''235.49'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '''5.14'
$ws.Range("E23").Value = '  -3.24%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = '''2.47'
$ws.Range("E25").Value = '  -10.33%  '
$ws.Range("D26").Value = '''2.26'
$ws.Range("E26").Value = '  -2.91%  '
$ws.Range("D27").Value = '''9.56'
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("D28").Value = '''157.79'
$ws.Range("E28").Value = '  -4.58%  '
$ws.Range("E29").Value = '  +20.19%  '
$ws.Range("D30").Value = '''19.52'
$ws.Range("E30").Value = '  -1.29%  '
$ws.Range("D31").Value = '''0.118'
$ws.Range("E31").Value = '  -2.19%  '
$ws.Range("D32").Value = '''4.80'
$ws.Range("E32").Value = '  -7.20%  '
$ws.Range("D33").Value = '''1.12'
$ws.Range("E33").Value = '  -7.04%  '
$ws.Range("D34").Value = '''0.0609'
$ws.Range("E34").Value = '  -0.54%  '
$ws.Range("D35").Value = '''4.30'
$ws.Range("E35").Value = '  -8.07%  '
$ws.Range("B36").Value = 'BinanceUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D36").Value = '''1.01'
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").Value = '''6.16'
$ws.Range("E37").Value = '  +3.59%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").Value = '''1.83'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''2.26'
$ws.Range("E39").Value = '  -6.92%  '
$ws.Range("D40").Value = '''3.06'
$ws.Range("E40").Value = '  +13.28%  '
$ws.Range("D41").Value = '''0.0972'
$ws.Range("E41").Value = '  -6.73%  '
$ws.Range("D42").Value = '''1.20'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("D43").Value = '''2.82'
$ws.Range("E43").Value = '  -4.08%  '
$ws.Range("D44").Value = '''0.0209'
$ws.Range("E44").Value = '  -4.02%  '
$ws.Range("D45").Value = '''1.07'
$ws.Range("E45").Value = '  -5.25%  '
$ws.Range("D46").Value = '''91.36'
$ws.Range("E46").Value = '  -3.54%  '
$ws.Range("D47").Value = '''15.79'
$ws.Range("E47").Value = '  -5.77%  '
$ws.Range("D48").Value = '''7.45'
$ws.Range("E48").Value = '  -7.18%  '
$ws.Range("D49").Value = '1.322.88'
$ws.Range("E49").Value = '  -6.53%  '
$ws.Range("D50").Value = '''2.72'
$ws.Range("E50").Value = '  -7.16%  '
$ws.Range("D51").Value = '2.152.94'
$ws.Range("E51").Value = '  -3.20%  '
